$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 70 (date + prices move from 45119 to 45147) ---
$ws.Cells.Item(70, 4).Value  = 45147   # D70 Fecha
$ws.Cells.Item(70, 14).Value = 14000   # N70 Precio minimo
$ws.Cells.Item(70, 15).Value = 14000   # O70 Precio maximo
$ws.Cells.Item(70, 16).Value = 14000   # P70 Precio promedio ponderado
$ws.Cells.Item(70, 19).Value = 778     # S70 Precio $/Kg

# --- Update existing row 71 (date + prices move from 45119 to 45147) ---
$ws.Cells.Item(71, 4).Value  = 45147   # D71 Fecha
$ws.Cells.Item(71, 14).Value = 12000   # N71 Precio minimo
$ws.Cells.Item(71, 15).Value = 12000   # O71 Precio maximo
$ws.Cells.Item(71, 16).Value = 12000   # P71 Precio promedio ponderado
$ws.Cells.Item(71, 19).Value = 667     # S71 Precio $/Kg

# --- Update existing row 72: becomes the old "Primera / 45119 / $/bandeja" record ---
$ws.Cells.Item(72, 4).Value  = 45119                          # D72 Fecha
$ws.Cells.Item(72, 12).Value = "Primera"                      # L72 Calidad
$ws.Cells.Item(72, 13).Value = 60                              # M72 Volumen
$ws.Cells.Item(72, 14).Value = 10000                           # N72 Precio minimo
$ws.Cells.Item(72, 15).Value = 10000                           # O72 Precio maximo
$ws.Cells.Item(72, 16).Value = 10000                           # P72 Precio promedio ponderado
$ws.Cells.Item(72, 17).Value = "$/bandeja 18 kilos granel"     # Q72 Unidad de comercializacion
$ws.Cells.Item(72, 19).Value = 556                              # S72 Precio $/Kg

# --- Update existing row 73: becomes the old "Segunda / 45119 / $/bandeja" record ---
$ws.Cells.Item(73, 4).Value  = 45119                          # D73 Fecha
$ws.Cells.Item(73, 12).Value = "Segunda"                      # L73 Calidad
$ws.Cells.Item(73, 13).Value = 60                              # M73 Volumen
$ws.Cells.Item(73, 14).Value = 8000                             # N73 Precio minimo
$ws.Cells.Item(73, 15).Value = 8000                             # O73 Precio maximo
$ws.Cells.Item(73, 16).Value = 8000                             # P73 Precio promedio ponderado
$ws.Cells.Item(73, 17).Value = "$/bandeja 18 kilos granel"     # Q73 Unidad de comercializacion
$ws.Cells.Item(73, 19).Value = 444                              # S73 Precio $/Kg

# --- New row 74: the old "Especial / 45043 / $/caja" record, shifted down ---
$ws.Cells.Item(74, 1).Value  = 7
$ws.Cells.Item(74, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(74, 3).Value  = "Ñuble"
$ws.Cells.Item(74, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(74, 4).Value  = 45043
$ws.Cells.Item(74, 5).Value  = 16
$ws.Cells.Item(74, 6).Value  = "Fruta"
$ws.Cells.Item(74, 7).Value  = 100104
$ws.Cells.Item(74, 8).Value  = "Frutos de pepita"
$ws.Cells.Item(74, 9).Value  = 100104003
$ws.Cells.Item(74, 10).Value = "Membrillo"
$ws.Cells.Item(74, 11).Value = "Champion"
$ws.Cells.Item(74, 12).Value = "Especial"
$ws.Cells.Item(74, 13).Value = 40
$ws.Cells.Item(74, 14).Value = 13000
$ws.Cells.Item(74, 15).Value = 13000
$ws.Cells.Item(74, 16).Value = 13000
$ws.Cells.Item(74, 17).Value = "$/caja 18 kilos empedrada"
$ws.Cells.Item(74, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(74, 19).Value = 722
$ws.Cells.Item(74, 20).Value = 18

# --- New row 75: the old "Primera / 45043 / $/caja" record, shifted down ---
$ws.Cells.Item(75, 1).Value  = 7
$ws.Cells.Item(75, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(75, 3).Value  = "Ñuble"
$ws.Cells.Item(75, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(75, 4).Value  = 45043
$ws.Cells.Item(75, 5).Value  = 16
$ws.Cells.Item(75, 6).Value  = "Fruta"
$ws.Cells.Item(75, 7).Value  = 100104
$ws.Cells.Item(75, 8).Value  = "Frutos de pepita"
$ws.Cells.Item(75, 9).Value  = 100104003
$ws.Cells.Item(75, 10).Value = "Membrillo"
$ws.Cells.Item(75, 11).Value = "Champion"
$ws.Cells.Item(75, 12).Value = "Primera"
$ws.Cells.Item(75, 13).Value = 50
$ws.Cells.Item(75, 14).Value = 12000
$ws.Cells.Item(75, 15).Value = 12000
$ws.Cells.Item(75, 16).Value = 12000
$ws.Cells.Item(75, 17).Value = "$/caja 18 kilos empedrada"
$ws.Cells.Item(75, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(75, 19).Value = 667
$ws.Cells.Item(75, 20).Value = 18
